$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-09-19"

# Update the label for the September row
$ws.Range("A10").Value = "September (through 09-19)"

# Update the September row values (row 10)
$ws.Range("B10").Value = 21
$ws.Range("C10").Value = 34
$ws.Range("D10").Value = 42
$ws.Range("E10").Value = 36
$ws.Range("F10").Value = 44
$ws.Range("G10").Value = 74
$ws.Range("H10").Value = 109

# Update the Total row values (row 11)
$ws.Range("B11").Value = 215
$ws.Range("C11").Value = 415
$ws.Range("D11").Value = 593
$ws.Range("E11").Value = 526
$ws.Range("F11").Value = 393
$ws.Range("G11").Value = 858
$ws.Range("H11").Value = 1179
